# Generate Report for Handoff
# Updates the status from "In Translation" to "Ready for handoff" and refreshes
# the handoff timestamps across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-31-21 00:31:09"

# zh-cn sheet: Status column C, Latest Handoff Datetime column E
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-21 00:31:06"

# de-de sheet: Status column C, Latest Handoff Datetime column E
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-21 00:31:09"
